$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.657.03"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "2.284.93"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.70%  "
$ws.Range("D5").Value = "'113.37"
$ws.Range("E5").Value = "  +19.37%  "
$ws.Range("D6").Value = "'268.07"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.622"
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").Value = "'0.614"
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("D10").Value = "'47.87"
$ws.Range("E10").Value = "  +7.02%  "
$ws.Range("D11").Value = "'0.0934"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").Value = "'8.59"
$ws.Range("E12").Value = "  +10.11%  "
$ws.Range("E13").Value = "  +1.93%  "
$ws.Range("D14").Value = "'15.51"
$ws.Range("E14").Value = "  +2.87%  "
$ws.Range("D15").Value = "2.639.75"
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").Value = "'0.849"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "2.299.74"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "43.773.91"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("D19").Value = "'0.0000109"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").Value = "'6.52"
$ws.Range("E20").Value = "  +5.32%  "
$ws.Range("D21").Value = "'72.26"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "'2.53"
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("D23").Value = "'232.33"
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("D24").Value = "'9.52"
$ws.Range("E24").Value = "  +4.97%  "
$ws.Range("D25").Value = "'2.83"
$ws.Range("E25").Value = "  +13.32%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'11.37"
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "'43.71"
$ws.Range("E28").Value = "  +9.28%  "
$ws.Range("D29").Value = "'3.42"
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("D31").Value = "'176.18"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'21.59"
$ws.Range("E32").Value = "  -3.24%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0925"
$ws.Range("E33").Value = "  +4.77%  "
$ws.Range("D34").Value = "'5.48"
$ws.Range("E34").Value = "  +2.24%  "
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").Value = "'4.69"
$ws.Range("E36").Value = "  +7.11%  "
$ws.Range("D37").Value = "'0.109"
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("D38").Value = "'3.88"
$ws.Range("E38").Value = "  +17.40%  "
$ws.Range("D39").Value = "'0.0353"
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("D40").Value = "'75.15"
$ws.Range("E40").Value = "  +15.45%  "
$ws.Range("D41").Value = "'0.241"
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("D42").Value = "'2.38"
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("D43").Value = "'13.13"
$ws.Range("E43").Value = "  +8.74%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'1.41"
$ws.Range("E45").Value = "  +5.75%  "
$ws.Range("D46").Value = "'5.88"
$ws.Range("E46").Value = "  +12.43%  "
$ws.Range("D47").Value = "'8.73"
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("D48").Value = "'0.100"
$ws.Range("E48").Value = "  -1.78%  "
$ws.Range("D49").Value = "'100.67"
$ws.Range("E49").Value = "  +2.52%  "
$ws.Range("D50").Value = "'1.23"
$ws.Range("E50").Value = "  +3.27%  "
$ws.Range("E51").Value = "  +6.08%  "
